# Create TestData for Invalid Login TestCase
# Adds a new "InvalidLogin" worksheet (after the existing "ValidLogin" sheet)
# containing a UserName/Password header row and an invalid abcd/xyz data row.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the ValidLogin sheet so it becomes sheet2
# (and, because it's newly created, Excel makes it the active/selected tab).
$validLogin = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $validLogin)
$ws.Name = "InvalidLogin"

# Header row
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Invalid credentials data row
$ws.Range("A2").Value = "abcd"
$ws.Range("B2").Value = "xyz"

# Match the column width tweak captured for column A on the new sheet.
$ws.Columns.Item(1).ColumnWidth = 11

# Leave the same cell selected/active on the new sheet as in the source workbook.
[void]$ws.Range("J16").Select()
